$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Sheet1")

# Update the text in C11 (shared string) from the old reminder text to the new one
$ws.Range("C11").Value = "Mieleenpalautus :D Home ja Nba komponentteja"

# Update hours value in B11 from 40 to 120 (this also causes I1 array formula to recalc)
$ws.Range("B11").Value = 120

# Update the active cell selection shown in the sheet view
$ws.Range("B12").Select()
